$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) of the last existing date cell (A464) so new date cells
# reuse the same cellXfs style index instead of creating a new one.
$formatSource = $ws.Range("A464")

$ws.Cells.Item(465,1).Value = 44539
$ws.Cells.Item(465,2).Value = 4
$ws.Cells.Item(465,3).Value = 16
$ws.Cells.Item(465,4).Value = 406.9175991861648

$ws.Cells.Item(466,1).Value = 44540
$ws.Cells.Item(466,2).Value = 3
$ws.Cells.Item(466,3).Value = 19
$ws.Cells.Item(466,4).Value = 483.2146490335707

$ws.Cells.Item(467,1).Value = 44541
$ws.Cells.Item(467,2).Value = 2
$ws.Cells.Item(467,3).Value = 21
$ws.Cells.Item(467,4).Value = 534.0793489318413

$ws.Cells.Item(468,1).Value = 44542
$ws.Cells.Item(468,2).Value = 1
$ws.Cells.Item(468,3).Value = 17
$ws.Cells.Item(468,4).Value = 432.3499491353002

$ws.Cells.Item(469,1).Value = 44543
$ws.Cells.Item(469,2).Value = 4
$ws.Cells.Item(469,3).Value = 20
$ws.Cells.Item(469,4).Value = 508.646998982706

$ws.Cells.Item(470,1).Value = 44544
$ws.Cells.Item(470,2).Value = 1
$ws.Cells.Item(470,3).Value = 16
$ws.Cells.Item(470,4).Value = 406.9175991861648

$ws.Cells.Item(471,1).Value = 44545
$ws.Cells.Item(471,2).Value = 0
$ws.Cells.Item(471,3).Value = 15
$ws.Cells.Item(471,4).Value = 381.4852492370295

$ws.Cells.Item(472,1).Value = 44546
$ws.Cells.Item(472,2).Value = 1
$ws.Cells.Item(472,3).Value = 12
$ws.Cells.Item(472,4).Value = 305.1881993896236

$ws.Cells.Item(473,1).Value = 44547
$ws.Cells.Item(473,2).Value = 7
$ws.Cells.Item(473,3).Value = 16
$ws.Cells.Item(473,4).Value = 406.9175991861648

$ws.Cells.Item(474,1).Value = 44548
$ws.Cells.Item(474,2).Value = 0
$ws.Cells.Item(474,3).Value = 14
$ws.Cells.Item(474,4).Value = 356.0528992878942

$ws.Cells.Item(475,1).Value = 44550
$ws.Cells.Item(475,2).Value = 1
$ws.Cells.Item(475,3).Value = 14
$ws.Cells.Item(475,4).Value = 356.0528992878942

$ws.Cells.Item(476,1).Value = 44551
$ws.Cells.Item(476,2).Value = 3
$ws.Cells.Item(476,3).Value = 13
$ws.Cells.Item(476,4).Value = 330.6205493387589

$ws.Cells.Item(477,1).Value = 44552
$ws.Cells.Item(477,2).Value = 1
$ws.Cells.Item(477,3).Value = 13
$ws.Cells.Item(477,4).Value = 330.6205493387589

$ws.Cells.Item(478,1).Value = 44553
$ws.Cells.Item(478,2).Value = 4
$ws.Cells.Item(478,3).Value = 17
$ws.Cells.Item(478,4).Value = 432.3499491353002

$ws.Cells.Item(479,1).Value = 44554
$ws.Cells.Item(479,2).Value = 4
$ws.Cells.Item(479,3).Value = 20
$ws.Cells.Item(479,4).Value = 508.646998982706

$ws.Cells.Item(480,1).Value = 44555
$ws.Cells.Item(480,2).Value = 3
$ws.Cells.Item(480,3).Value = 16
$ws.Cells.Item(480,4).Value = 406.9175991861648

$ws.Cells.Item(481,1).Value = 44556
$ws.Cells.Item(481,2).Value = 5
$ws.Cells.Item(481,3).Value = 21
$ws.Cells.Item(481,4).Value = 534.0793489318413

$ws.Cells.Item(482,1).Value = 44557
$ws.Cells.Item(482,2).Value = 1
$ws.Cells.Item(482,3).Value = 21
$ws.Cells.Item(482,4).Value = 534.0793489318413

$ws.Cells.Item(483,1).Value = 44558
$ws.Cells.Item(483,2).Value = 0
$ws.Cells.Item(483,3).Value = 18
$ws.Cells.Item(483,4).Value = 457.7822990844354

$ws.Cells.Item(484,1).Value = 44559
$ws.Cells.Item(484,2).Value = 0
$ws.Cells.Item(484,3).Value = 17
$ws.Cells.Item(484,4).Value = 432.3499491353002

$ws.Cells.Item(485,1).Value = 44560
$ws.Cells.Item(485,2).Value = 3
$ws.Cells.Item(485,3).Value = 16
$ws.Cells.Item(485,4).Value = 406.9175991861648

$ws.Cells.Item(486,1).Value = 44561
$ws.Cells.Item(486,2).Value = 6
$ws.Cells.Item(486,3).Value = 18
$ws.Cells.Item(486,4).Value = 457.7822990844354

$ws.Cells.Item(487,1).Value = 44562
$ws.Cells.Item(487,2).Value = 5
$ws.Cells.Item(487,3).Value = 20
$ws.Cells.Item(487,4).Value = 508.646998982706

$ws.Cells.Item(488,1).Value = 44563
$ws.Cells.Item(488,2).Value = 3
$ws.Cells.Item(488,3).Value = 18
$ws.Cells.Item(488,4).Value = 457.7822990844354

$ws.Cells.Item(489,1).Value = 44564
$ws.Cells.Item(489,2).Value = 2
$ws.Cells.Item(489,3).Value = 19
$ws.Cells.Item(489,4).Value = 483.2146490335707

$ws.Cells.Item(490,1).Value = 44565
$ws.Cells.Item(490,2).Value = 1
$ws.Cells.Item(490,3).Value = 20
$ws.Cells.Item(490,4).Value = 508.646998982706

$ws.Cells.Item(491,1).Value = 44566
$ws.Cells.Item(491,2).Value = 3
$ws.Cells.Item(491,3).Value = 23
$ws.Cells.Item(491,4).Value = 584.9440488301119

# Apply the same number format / style as column A uses elsewhere to the newly
# added date cells (A465:A491) by copying formats from A464.
$formatSource.Copy()
$ws.Range("A465:A491").PasteSpecial(-4122)
$excel.CutCopyMode = 0
